$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellXml($cell, $bodyXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $cell.Range.InsertXML($xml)
}

# Row 7: "Autosave of cmd file..." -> "Completed and tested"; this is also where the
# "_GoBack" bookmark (last-edit marker) now lives, so it is added here.
Set-CellXml $t.Cell(7, 2) '<w:p><w:r><w:t>Completed and tested</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Row 8: "Fill empty bookmarks..." -> "Next Rev"
$t.Cell(8, 2).Range.Text = "Next Rev"

# Row 10: "Page numbering and footnotes look odd on landscape pages." -> "No issue found"
$t.Cell(10, 2).Range.Text = "No issue found"

# Row 11: "Automatically remove blank pages" -> "Not implementing for now"
$t.Cell(11, 2).Range.Text = "Not implementing for now"

# Row 12: "Ability to re-order pages..." -> "Not implementing for now"
$t.Cell(12, 2).Range.Text = "Not implementing for now"

# Row 13: "Sort columns by name or date" -> "Not implementing for now"
$t.Cell(13, 2).Range.Text = "Not implementing for now"

# Row 14: "Ability to prefix and restart page numbering per file..." -> "Next Rev"; the
# "_GoBack" bookmark that used to sit in this (previously empty) cell moves to row 7 above,
# so it is dropped from here.
Set-CellXml $t.Cell(14, 2) '<w:p><w:r><w:t>Next Rev</w:t></w:r></w:p>'
